# "Primeira versão do jogo"
#
# Applies the authored changes to "planejamento - IGA.xlsx" (sheet DATAS - IGA):
#   1. Cell E4: reword the task text to mention "frameworks..." as well.
#   2. Cell B9: responsible person corrected from "Heitor/Daniel" to "Heitor".
#   3. View state: zoom 130% -> 145%, scrolled so row 2 is at the top, and the
#      active selection moved to E4.
#   4. The 9 "Forms" checkboxes anchored in column C (rows 1-9) shift a hair to
#      the right (Excel re-snaps their legacy VML anchors to the pixel grid).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Content edits -----------------------------------------------------

$ws.Range("E4").Value = "Escolher linguagem e ferramentas (IDE, bibliotecas, frameworks...)"
$ws.Range("B9").Value = "Heitor"

# --- 2. Checkbox anchor nudge ----------------------------------------------
# All nine Form checkboxes sit in column C; only their horizontal (Left)
# position shifts - top/height are untouched. Compute the new Left from the
# live position of column C plus the exact target offset (EMU / 12700 = pt).

$colLeft = $ws.Range("C1").Left

$checkboxColOffsets = @{
    "Check Box 7"  = 180975
    "Check Box 8"  = 180975
    "Check Box 15" = 180975
    "Check Box 17" = 190500
    "Check Box 18" = 180975
    "Check Box 19" = 190500
    "Check Box 21" = 180975
    "Check Box 22" = 190500
    "Check Box 23" = 180975
}

foreach ($name in $checkboxColOffsets.Keys) {
    $targetColOffEmu = $checkboxColOffsets[$name]
    $shp = $ws.Shapes.Item($name)
    $shp.Left = $colLeft + ($targetColOffEmu / 12700.0)
}

# --- 3. View state ----------------------------------------------------------

$win = $excel.ActiveWindow
$win.Zoom = 145
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("E4").Select()
